# "Implementado Dublin Core en tablas"
# - Remove the old "bga-obra" summary sheet (kept only as a draft/merge of
#   the "obra" and "referente" tables).
# - Re-label the headers of the "obra" and "referente" tables using
#   Dublin Core terms (Title, Date, Medium, Format, Publisher) instead of
#   the old Spanish labels.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# 1. Drop the "bga-obra" sheet entirely.
$bga = $wb.Worksheets.Item("bga-obra")
[void]$bga.Delete()

# 2. "obra" table -> Dublin Core headers.
#    ID | Título | Fecha | Dimensiones | Técnica | Archivo
# -> ID | Title  | Date  | Format      | Medium  | Archivo
$obra = $wb.Worksheets.Item("obra")
$obra.Range("B1").Value = "Title"
$obra.Range("C1").Value = "Date"
$obra.Range("E1").Value = "Medium"
$obra.Range("D1").Value = "Format"

# 3. "referente" table -> Dublin Core headers.
#    ID | Título | Fecha | Periódico  | Archivo
# -> ID | Title  | Date  | Publisher  | Archivo
$referente = $wb.Worksheets.Item("referente")
$referente.Range("B1").Value = "Title"
$referente.Range("C1").Value = "Date"
$referente.Range("D1").Value = "Publisher"

# 4. Restore cursor positions / active sheet, matching the saved state.
[void]$obra.Activate()
[void]$obra.Range("D9").Select()

[void]$referente.Range("D13").Select()

[void]$obra.Activate()
